# Kareli.xlsx – "upgrade left table until javakheti"
#
# 1) Rename the sheet tab from "1" to "Kareli"
# 2) Mask a handful of now-unreliable observations in the Urban/Rural rows
#    with the confidentiality markers already used elsewhere in the table
#    ("..." for brand-new suppressions, "…" for the ones that reuse the
#    existing note glyph)
# 3) Remove the blank spacer row between the data block and the footnote,
#    so the footnote moves up from row 9 to row 8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Kareli"

# Row 6 = Urban
$ws.Range("B6").Value = "..."
$ws.Range("C6").Value = "..."
$ws.Range("H6").Value = "..."
$ws.Range("I6").Value = "..."
$ws.Range("J6").Value = "..."
$ws.Range("K6").Value = "…"
$ws.Range("M6").Value = "…"

# Row 7 = Rural
$ws.Range("B7").Value = "..."
$ws.Range("C7").Value = "..."
$ws.Range("H7").Value = "..."
$ws.Range("K7").Value = "…"

# Drop the empty row 8 so the footnote row (was 9) becomes row 8
$ws.Rows.Item(8).Delete()
